$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

# Mapping of row label (column A) -> new hashcode value (column B)
# Generated from the commit's automatic hashcode update diff.
$updates = @(
    @{ Row = 9; Label = "05-050305TC"; OldHash = "c76eb7bc1a2e6e67034fc7d750762de9"; NewHash = "09768edd95a8b219f10218dc50a94417" }
    @{ Row = 17; Label = "05-050305TP"; OldHash = "6d8ffd1d66c53c710be8772851e1d28a"; NewHash = "1ccfc1ec97dfed9f35c1ed5011b1cea9" }
    @{ Row = 126; Label = "05-050309A"; OldHash = "885125c12218fb55a9c17a473ea811f5"; NewHash = "f8e2c2e76e50c47fd884009976743833" }
    @{ Row = 136; Label = "05-050312TC"; OldHash = "b7039b57dbda92005e340256ad999b90"; NewHash = "eb7b0979e989c558249db2170fe6a48d" }
    @{ Row = 159; Label = "05-050203TP"; OldHash = "e180276079263c04640119ac2f9a2356"; NewHash = "4749c882ce4f82f5ec89fee91ecc415c" }
    @{ Row = 169; Label = "05-050203TC"; OldHash = "57c8ebb0b1bfe05484cfbeee6e543676"; NewHash = "4da83de0fa8baa0c3e34ef948fa497bf" }
    @{ Row = 183; Label = "05-050305A"; OldHash = "1566ad624c9b683444f8640e7090cefd"; NewHash = "477b146f8b21754abe9e6418d07f97ae" }
    @{ Row = 200; Label = "05-050306A"; OldHash = "d5ef55e36803ff9c65c83cdd13fffe52"; NewHash = "875decfdb4d3f6746c65a89f45459306" }
    @{ Row = 228; Label = "05-050304A"; OldHash = "fe38701a3da4b84079059572acfcc9b3"; NewHash = "5b813c348de89f8832b3df7554abeb70" }
    @{ Row = 281; Label = "05-050201TC"; OldHash = "5303e7c7c414586e96e97fca9adc5a1a"; NewHash = "afc91a4d0896544a39504d970bebe301" }
    @{ Row = 302; Label = "05-050310TP"; OldHash = "c72ab92478c61d71a94c691b800f69f1"; NewHash = "41f7a08e5604f7733de62b092e819c2d" }
    @{ Row = 339; Label = "05-050201TP"; OldHash = "099ff95134ac2a6dda1c6112387b1c53"; NewHash = "1eb832b6afed5fa4baf694d891211e50" }
    @{ Row = 464; Label = "05-050204A"; OldHash = "88ca15026fa327f90edcf2607339c165"; NewHash = "3c75af0a389448ba653dbb96b057f85d" }
    @{ Row = 506; Label = "05-050202A"; OldHash = "d11e0cd41977733cd7b40226af342944"; NewHash = "4670f7f253d8abe8a660119fd708e885" }
    @{ Row = 507; Label = "05-050311A"; OldHash = "28b7f4082aa807fa960d3091d6953006"; NewHash = "bcf10a301975099317a3671d48f56727" }
    @{ Row = 524; Label = "05-050203A"; OldHash = "b442e64966200cb4be835787721f9bae"; NewHash = "47b1b203b6ab8a70b7b10583d0108c5b" }
    @{ Row = 558; Label = "05-050310A"; OldHash = "cb211322d39ea5dcae043e1ec1002c9b"; NewHash = "500fec36363758d7e706ee1f3a320cbd" }
    @{ Row = 580; Label = "05-050308TP"; OldHash = "521ce29e8304ca26acab34907e3d08da"; NewHash = "90e9978e5fac4cdc1c413f6cc4049a3c" }
    @{ Row = 624; Label = "05-050204TP"; OldHash = "8eed330081db7ea415c2ac50c2458014"; NewHash = "a619418188285d32ee4afa2a1af3c1ad" }
    @{ Row = 635; Label = "05-050204TC"; OldHash = "d450c3da6f90944d2dbd85eeeee6c17e"; NewHash = "eff5797203762a41ac372a1640233c11" }
    @{ Row = 688; Label = "05-050206TP"; OldHash = "87f7d8c8d5f14748512c9245c79f6ea6"; NewHash = "7d2715d459ca0f0c3f692fd10702b608" }
    @{ Row = 693; Label = "05-050206TC"; OldHash = "e992428de39ad6cc52cb72f089587295"; NewHash = "ab23ac1348387edfb9f7c498fb3e5f2a" }
    @{ Row = 708; Label = "05-050304TC"; OldHash = "12e5dbeb119384264be0298d3ffb04dd"; NewHash = "c73244e4d02da93b2f5418460dd36c9d" }
    @{ Row = 711; Label = "05-050206A"; OldHash = "d174fa8fbca0c777f41402c2571309ad"; NewHash = "85376c330cb8c179172eb8012e4289fc" }
    @{ Row = 723; Label = "05-050304TP"; OldHash = "3d55dde6eea0e77c61e852a4347905de"; NewHash = "c5ee5882e46f01af84add9b219ddf0c2" }
    @{ Row = 827; Label = "05-050202TP"; OldHash = "6f14a86add7ba4c658e6672d743c2b75"; NewHash = "3cad1c31d6cda35f1ce8b17cbb9cfdb9" }
    @{ Row = 838; Label = "05-050311TC"; OldHash = "10e0d3fcba82c94ccc94802d6c5c9179"; NewHash = "e08d817cc6a46610a3b5f893585aa94e" }
    @{ Row = 843; Label = "05-050311TP"; OldHash = "08ec81e9257330f99b6ec686fc7b6d56"; NewHash = "b102e7c044aa28ec0c96f4f071d794ab" }
    @{ Row = 870; Label = "05-050309TP"; OldHash = "f5701873862730c7a2f060b7921941ab"; NewHash = "376b400271a9aac22e19182e385681ae" }
)

foreach ($u in $updates) {
    $labelCell = $ws.Cells.Item($u.Row, 1)
    $hashCell  = $ws.Cells.Item($u.Row, 2)
    if ($labelCell.Value -eq $u.Label -and $hashCell.Value -eq $u.OldHash) {
        $hashCell.Value = $u.NewHash
    } else {
        # Fallback: search the used range for the label in column A in case rows shifted
        $found = $ws.Columns.Item(1).Find($u.Label, [System.Reflection.Missing]::Value, -4163, 1)
        if ($found -ne $null) {
            $ws.Cells.Item($found.Row, 2).Value = $u.NewHash
        }
    }
}
